# Update the MAKE for plate ABC123 (row 7) from FORD to JAGUAR.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "JAGUAR"

# Match the author's final selection as captured in the diff (B7 selected).
$ws.Range("B7").Select()
